$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting from the cell above (C9) onto the new cell (C10) so the
# newly-populated "Company code" cell for velpatasvir picks up the same
# style index as its neighbours, then set its value to the new shared string.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C10").Value = "GS-5816"

# Reflect the widened selection (now that column C's data extends through
# row 12) in the sheet view.
$ws.Range("A1:C12").Select()
